# 02-Cleaned_Dataset.xlsx edit script
# 1) Fix mojibake: double every literal "EE" substring that appears in tweet text
#    (case-sensitive partial match) -> "EEEE"
# 2) Normalize/clean up emotion labels in columns (trim trailing spaces, fold
#    near-duplicate labels into the canonical spelling), which also removes the
#    now-unused duplicate shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# xlPart = 2, xlWhole = 1
$xlPart = 2
$xlWhole = 1

# Step 1: double up "EE" -> "EEEE" wherever it occurs (case sensitive, partial match)
$used.Replace("EE", "EEEE", $xlPart, $xlWhole, $true, $false, $false, $false) | Out-Null

# Step 2: normalize emotion-label cells (exact, case-sensitive whole-cell matches)
$used.Replace("guilt", "guilty", $xlWhole, $xlWhole, $true, $false, $false, $false) | Out-Null
$used.Replace("lonely ", "lonely", $xlWhole, $xlWhole, $true, $false, $false, $false) | Out-Null
$used.Replace("depressed ", "depressed", $xlWhole, $xlWhole, $true, $false, $false, $false) | Out-Null
$used.Replace("displeased ", "displeased", $xlWhole, $xlWhole, $true, $false, $false, $false) | Out-Null
$used.Replace("grief ", "grief", $xlWhole, $xlWhole, $true, $false, $false, $false) | Out-Null
$used.Replace("guilt ", "guilty", $xlWhole, $xlWhole, $true, $false, $false, $false) | Out-Null
$used.Replace("lost", "lonely", $xlWhole, $xlWhole, $true, $false, $false, $false) | Out-Null
$used.Replace("depression", "depressed", $xlWhole, $xlWhole, $true, $false, $false, $false) | Out-Null
